$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Junio de 2020 a las 19:58"

# Swap adjacent country-name pairs whose ranking crossed after the refresh
$ws.Range("A102").Value = "Maldivas"
$ws.Range("A103").Value = "Guayana Francesa"
$ws.Range("A137").Value = "Estado de Palestina"
$ws.Range("A138").Value = "Uganda"
$ws.Range("A162").Value = "Comoras"
$ws.Range("A163").Value = "Martinica"
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("A209").Value = "Santa Sede"
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("A214").Value = "Papua Nueva Guinea"

# Refresh case-count statistics (Casos totales / Nuevos casos / Casos activos / Recuperados / Casos criticos / Muertes)
$ws.Range("B4").Value = 2314549
$ws.Range("C4").Value = 17359
$ws.Range("D4").Value = 958234
$ws.Range("E4").Value = 1234641
$ws.Range("G4").Value = 267
$ws.Range("H4").Value = 121674
$ws.Range("B5").Value = 1043168
$ws.Range("C5").Value = 4600
$ws.Range("D5").Value = 543186
$ws.Range("E5").Value = 450826
$ws.Range("G5").Value = 66
$ws.Range("H5").Value = 49156
$ws.Range("B7").Value = 407689
$ws.Range("C7").Value = 11877
$ws.Range("D7").Value = 220349
$ws.Range("E7").Value = 174071
$ws.Range("G7").Value = 299
$ws.Range("H7").Value = 13269
$ws.Range("D12").Value = 196609
$ws.Range("E12").Value = 35844
$ws.Range("B14").Value = 190965
$ws.Range("C14").Value = 305
$ws.Range("E14").Value = 7605
$ws.Range("B15").Value = 186493
$ws.Range("C15").Value = 1248
$ws.Range("D15").Value = 158828
$ws.Range("E15").Value = 22738
$ws.Range("G15").Value = 22
$ws.Range("H15").Value = 4927
$ws.Range("B21").Value = 100959
$ws.Range("C21").Value = 330
$ws.Range("D21").Value = 63450
$ws.Range("E21").Value = 29099
$ws.Range("G21").Value = 64
$ws.Range("H21").Value = 8410
$ws.Range("B46").Value = 25374
$ws.Range("C46").Value = 6
$ws.Range("E46").Value = 961
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 1715
$ws.Range("E56").Value = 5990
$ws.Range("G56").Value = 5
$ws.Range("H56").Value = 118
$ws.Range("B68").Value = 9839
$ws.Range("C68").Value = 226
$ws.Range("D68").Value = 8223
$ws.Range("E68").Value = 1403
$ws.Range("B96").Value = 2755
$ws.Range("C96").Value = 36
$ws.Range("D96").Value = 751
$ws.Range("E96").Value = 1916
$ws.Range("B102").Value = 2187
$ws.Range("C102").Value = 37
$ws.Range("D102").Value = 1788
$ws.Range("E102").Value = 391
$ws.Range("H102").Value = 8
$ws.Range("B103").Value = 2163
$ws.Range("C103").Value = 194
$ws.Range("D103").Value = 890
$ws.Range("E103").Value = 1268
$ws.Range("H103").Value = 5
$ws.Range("B129").Value = 922
$ws.Range("C129").Value = 3
$ws.Range("D129").Value = 328
$ws.Range("E129").Value = 340
$ws.Range("G129").Value = 3
$ws.Range("H129").Value = 254
$ws.Range("B137").Value = 785
$ws.Range("C137").Value = 110
$ws.Range("D137").Value = 437
$ws.Range("E137").Value = 345
$ws.Range("H137").Value = 3
$ws.Range("B138").Value = 763
$ws.Range("C138").Value = 8
$ws.Range("D138").Value = 492
$ws.Range("E138").Value = 271
$ws.Range("H138").Value = 0
$ws.Range("B153").Value = 505
$ws.Range("C153").Value = 1
$ws.Range("E153").Value = 44
$ws.Range("B156").Value = 359
$ws.Range("C156").Value = 4
$ws.Range("E156").Value = 35
$ws.Range("B162").Value = 247
$ws.Range("C162").Value = 37
$ws.Range("D162").Value = 159
$ws.Range("E162").Value = 83
$ws.Range("H162").Value = 5
$ws.Range("B163").Value = 236
$ws.Range("D163").Value = 98
$ws.Range("E163").Value = 124
$ws.Range("H163").Value = 14
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
